$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Read current data rows (16-31) into memory, then re-sort them by
# Periodo Mora (column E) ascending, then Tipo Doc Trabajador (column B)
# ascending, and write the values back - interleaving the two workers'
# records by period instead of grouping them by worker.
$startRow = 16
$endRow = 31

$docRank = @{ CC = 0; CE = 1 }

$data = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $row = @{
        B = $bVal
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $eVal
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        SortKey = ([int]$eVal * 10) + $docRank[$bVal]
    }
    $data += ,$row
}

# Sort ascending by period (E) first, then by doc type (B: CC before CE),
# via a single composite key computed up front.
$sorted = $data | Sort-Object -Property SortKey

$i = $startRow
foreach ($row in $sorted) {
    $ws.Cells.Item($i, 2).Value2 = $row.B
    $ws.Cells.Item($i, 3).Value2 = $row.C
    $ws.Cells.Item($i, 4).Value2 = $row.D
    $ws.Cells.Item($i, 5).Value2 = $row.E
    $ws.Cells.Item($i, 6).Value2 = $row.F
    $ws.Cells.Item($i, 7).Value2 = $row.G
    $i++
}

$wb.Save()
